# Fruta / hortaliza, semanal
# Insert a new weekly record into the "Granada" (Vega Modelo de Temuco) price
# series. The new observation is inserted as row 232, pushing the existing
# rows 232:268 down to 233:269 (dimension grows from A1:T268 to A1:T269).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 232 - Excel shifts rows 232:268
# down to 233:269 and carries the row's number formatting (date style on
# column D) down with them / onto the new row.
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new weekly observation.
$ws.Range("A232").Value2 = 10
$ws.Range("B232").Value2 = "Vega Modelo de Temuco"
$ws.Range("C232").Value2 = "La Araucanía"
$ws.Range("D232").Value2 = 45127
$ws.Range("E232").Value2 = 9
$ws.Range("F232").Value2 = "Fruta"
$ws.Range("G232").Value2 = 100104
$ws.Range("H232").Value2 = "Frutos de pepita"
$ws.Range("I232").Value2 = 100104001
$ws.Range("J232").Value2 = "Granada"
$ws.Range("K232").Value2 = "Wonderfull"
$ws.Range("L232").Value2 = "Primera"
$ws.Range("M232").Value2 = 35
$ws.Range("N232").Value2 = 16000
$ws.Range("O232").Value2 = 16000
$ws.Range("P232").Value2 = 16000
$ws.Range("Q232").Value2 = "$/bandeja 10 kilos granel"
$ws.Range("R232").Value2 = "Provincia de Limarí"
$ws.Range("S232").Value2 = 1600
$ws.Range("T232").Value2 = 10
